$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (label "0.0")
$ws.Range("B2").Value = 0.9646133682830931
$ws.Range("C2").Value = 0.9496774193548387
$ws.Range("D2").Value = 0.9570871261378414
$ws.Range("E2").Value = 1550

# Row 3 (label "1.0")
$ws.Range("B3").Value = 0.9470827679782904
$ws.Range("C3").Value = 0.9627586206896551
$ws.Range("D3").Value = 0.9548563611491109
$ws.Range("E3").Value = 1450

# Row 4 (label "accuracy")
$ws.Range("B4").Value = 0.956
$ws.Range("C4").Value = 0.956
$ws.Range("D4").Value = 0.956
$ws.Range("E4").Value = 0.956

# Row 5 (label "macro avg")
$ws.Range("B5").Value = 0.9558480681306918
$ws.Range("C5").Value = 0.956218020022247
$ws.Range("D5").Value = 0.9559717436434761
$ws.Range("E5").Value = 3000

# Row 6 (label "weighted avg")
$ws.Range("B6").Value = 0.9561402448024385
$ws.Range("C6").Value = 0.956
$ws.Range("D6").Value = 0.956008923059955
$ws.Range("E6").Value = 3000
